$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Update the values in E16:E22 from 5 to 10
$ws.Range("E16:E22").Value = 10

# Make DBD the active sheet and move the selection to G20
$ws.Activate()
$ws.Range("G20").Select()
